$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Risk probability/severity matrix (top grid) updates ---
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 7
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 8
$ws.Range("D13").Value = 6
$ws.Range("F13").Value = "4, 8"

# --- Existing risk table updates (item 6 probability/severity-after-mitigation) ---
$ws.Range("D24").Value = 3
$ws.Range("G24").Value = 2

# --- New risk item 8 row ---
$ws.Range("B26").Value = 8
$ws.Range("C26").Value = "Applikasjonen blir ikke det som oppdragsgiver hadde håpet på"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = "Møter annenhver uke for å diskutere å planlegge behov for brukere"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 3
